$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting (date number format, styles) from the last existing
# data row (2589) down across the new rows (2590:2629) in one shot,
# then overwrite each cell with the new data below.
$ws.Range("A2589:H2589").Copy($ws.Range("A2590:H2629")) | Out-Null

$ws.Cells.Item(2590, 1).Value = 44169
$ws.Cells.Item(2590, 2).Value = '0-10 years'
$ws.Cells.Item(2590, 3).Value = 19543
$ws.Cells.Item(2590, 4).Value = 0.0497773860950363
$ws.Cells.Item(2590, 5).Value = 211
$ws.Cells.Item(2590, 6).Value = 0.0484389348025712
$ws.Cells.Item(2590, 7).Value = 4
$ws.Cells.Item(2590, 8).Value = -1

$ws.Cells.Item(2591, 1).Value = 44169
$ws.Cells.Item(2591, 2).Value = '11-20 years'
$ws.Cells.Item(2591, 3).Value = 51805
$ws.Cells.Item(2591, 4).Value = 0.131950953622952
$ws.Cells.Item(2591, 5).Value = 475
$ws.Cells.Item(2591, 6).Value = 0.109044995408632
$ws.Cells.Item(2591, 7).Value = 2
$ws.Cells.Item(2591, 8).Value = -1

$ws.Cells.Item(2592, 1).Value = 44169
$ws.Cells.Item(2592, 2).Value = '21-30 years'
$ws.Cells.Item(2592, 3).Value = 74767
$ws.Cells.Item(2592, 4).Value = 0.190436771538023
$ws.Cells.Item(2592, 5).Value = 740
$ws.Cells.Item(2592, 6).Value = 0.169880624426079
$ws.Cells.Item(2592, 7).Value = 31
$ws.Cells.Item(2592, 8).Value = 0

$ws.Cells.Item(2593, 1).Value = 44169
$ws.Cells.Item(2593, 2).Value = '31-40 years'
$ws.Cells.Item(2593, 3).Value = 61715
$ws.Cells.Item(2593, 4).Value = 0.157192415844812
$ws.Cells.Item(2593, 5).Value = 700
$ws.Cells.Item(2593, 6).Value = 0.160697887970615
$ws.Cells.Item(2593, 7).Value = 62
$ws.Cells.Item(2593, 8).Value = 1

$ws.Cells.Item(2594, 1).Value = 44169
$ws.Cells.Item(2594, 2).Value = '41-50 years'
$ws.Cells.Item(2594, 3).Value = 58521
$ws.Cells.Item(2594, 4).Value = 0.149057074741218
$ws.Cells.Item(2594, 5).Value = 700
$ws.Cells.Item(2594, 6).Value = 0.160697887970615
$ws.Cells.Item(2594, 7).Value = 162
$ws.Cells.Item(2594, 8).Value = 1

$ws.Cells.Item(2595, 1).Value = 44169
$ws.Cells.Item(2595, 2).Value = '51-60 years'
$ws.Cells.Item(2595, 3).Value = 53330
$ws.Cells.Item(2595, 4).Value = 0.135835235145489
$ws.Cells.Item(2595, 5).Value = 665
$ws.Cells.Item(2595, 6).Value = 0.152662993572084
$ws.Cells.Item(2595, 7).Value = 450
$ws.Cells.Item(2595, 8).Value = 11

$ws.Cells.Item(2596, 1).Value = 44169
$ws.Cells.Item(2596, 2).Value = '61-70 years'
$ws.Cells.Item(2596, 3).Value = 37698
$ws.Cells.Item(2596, 4).Value = 0.0960194392371016
$ws.Cells.Item(2596, 5).Value = 462
$ws.Cells.Item(2596, 6).Value = 0.106060606060606
$ws.Cells.Item(2596, 7).Value = 912
$ws.Cells.Item(2596, 8).Value = 14

$ws.Cells.Item(2597, 1).Value = 44169
$ws.Cells.Item(2597, 2).Value = '71-80 years'
$ws.Cells.Item(2597, 3).Value = 22658
$ws.Cells.Item(2597, 4).Value = 0.0577115086804141
$ws.Cells.Item(2597, 5).Value = 268
$ws.Cells.Item(2597, 6).Value = 0.061524334251607
$ws.Cells.Item(2597, 7).Value = 1475
$ws.Cells.Item(2597, 8).Value = 28

$ws.Cells.Item(2598, 1).Value = 44169
$ws.Cells.Item(2598, 2).Value = '81+ years'
$ws.Cells.Item(2598, 3).Value = 12037
$ws.Cells.Item(2598, 4).Value = 0.0306590797946043
$ws.Cells.Item(2598, 5).Value = 148
$ws.Cells.Item(2598, 6).Value = 0.0339761248852158
$ws.Cells.Item(2598, 7).Value = 1778
$ws.Cells.Item(2598, 8).Value = 43

$ws.Cells.Item(2599, 1).Value = 44169
$ws.Cells.Item(2599, 2).Value = 'Pending'
$ws.Cells.Item(2599, 3).Value = 534
$ws.Cells.Item(2599, 4).Value = 0.00136013530035048
$ws.Cells.Item(2599, 5).Value = -13
$ws.Cells.Item(2599, 6).Value = -0.00298438934802571
$ws.Cells.Item(2599, 7).Value = 0
$ws.Cells.Item(2599, 8).Value = -1

$ws.Cells.Item(2600, 1).Value = 44170
$ws.Cells.Item(2600, 2).Value = '0-10 years'
$ws.Cells.Item(2600, 3).Value = 19773
$ws.Cells.Item(2600, 4).Value = 0.0497406432851515
$ws.Cells.Item(2600, 5).Value = 230
$ws.Cells.Item(2600, 6).Value = 0.0468050468050468
$ws.Cells.Item(2600, 7).Value = 4
$ws.Cells.Item(2600, 8).Value = 0

$ws.Cells.Item(2601, 1).Value = 44170
$ws.Cells.Item(2601, 2).Value = '11-20 years'
$ws.Cells.Item(2601, 3).Value = 52385
$ws.Cells.Item(2601, 4).Value = 0.131778870100271
$ws.Cells.Item(2601, 5).Value = 580
$ws.Cells.Item(2601, 6).Value = 0.118030118030118
$ws.Cells.Item(2601, 7).Value = 2
$ws.Cells.Item(2601, 8).Value = 0

$ws.Cells.Item(2602, 1).Value = 44170
$ws.Cells.Item(2602, 2).Value = '21-30 years'
$ws.Cells.Item(2602, 3).Value = 75621
$ws.Cells.Item(2602, 4).Value = 0.190230980926842
$ws.Cells.Item(2602, 5).Value = 854
$ws.Cells.Item(2602, 6).Value = 0.173789173789174
$ws.Cells.Item(2602, 7).Value = 31
$ws.Cells.Item(2602, 8).Value = 0

$ws.Cells.Item(2603, 1).Value = 44170
$ws.Cells.Item(2603, 2).Value = '31-40 years'
$ws.Cells.Item(2603, 3).Value = 62495
$ws.Cells.Item(2603, 4).Value = 0.15721142477649
$ws.Cells.Item(2603, 5).Value = 780
$ws.Cells.Item(2603, 6).Value = 0.158730158730159
$ws.Cells.Item(2603, 7).Value = 62
$ws.Cells.Item(2603, 8).Value = 0

$ws.Cells.Item(2604, 1).Value = 44170
$ws.Cells.Item(2604, 2).Value = '41-50 years'
$ws.Cells.Item(2604, 3).Value = 59292
$ws.Cells.Item(2604, 4).Value = 0.14915400908629
$ws.Cells.Item(2604, 5).Value = 771
$ws.Cells.Item(2604, 6).Value = 0.156898656898657
$ws.Cells.Item(2604, 7).Value = 162
$ws.Cells.Item(2604, 8).Value = 0

$ws.Cells.Item(2605, 1).Value = 44170
$ws.Cells.Item(2605, 2).Value = '51-60 years'
$ws.Cells.Item(2605, 3).Value = 54049
$ws.Cells.Item(2605, 4).Value = 0.135964801948068
$ws.Cells.Item(2605, 5).Value = 719
$ws.Cells.Item(2605, 6).Value = 0.146316646316646
$ws.Cells.Item(2605, 7).Value = 451
$ws.Cells.Item(2605, 8).Value = 1

$ws.Cells.Item(2606, 1).Value = 44170
$ws.Cells.Item(2606, 2).Value = '61-70 years'
$ws.Cells.Item(2606, 3).Value = 38175
$ws.Cells.Item(2606, 4).Value = 0.0960324208471481
$ws.Cells.Item(2606, 5).Value = 477
$ws.Cells.Item(2606, 6).Value = 0.0970695970695971
$ws.Cells.Item(2606, 7).Value = 916
$ws.Cells.Item(2606, 8).Value = 4

$ws.Cells.Item(2607, 1).Value = 44170
$ws.Cells.Item(2607, 2).Value = '71-80 years'
$ws.Cells.Item(2607, 3).Value = 22949
$ws.Cells.Item(2607, 4).Value = 0.0577301382061873
$ws.Cells.Item(2607, 5).Value = 291
$ws.Cells.Item(2607, 6).Value = 0.0592185592185592
$ws.Cells.Item(2607, 7).Value = 1484
$ws.Cells.Item(2607, 8).Value = 9

$ws.Cells.Item(2608, 1).Value = 44170
$ws.Cells.Item(2608, 2).Value = '81+ years'
$ws.Cells.Item(2608, 3).Value = 12236
$ws.Cells.Item(2608, 4).Value = 0.0307806863519503
$ws.Cells.Item(2608, 5).Value = 199
$ws.Cells.Item(2608, 6).Value = 0.0404965404965405
$ws.Cells.Item(2608, 7).Value = 1793
$ws.Cells.Item(2608, 8).Value = 15

$ws.Cells.Item(2609, 1).Value = 44170
$ws.Cells.Item(2609, 2).Value = 'Pending'
$ws.Cells.Item(2609, 3).Value = 547
$ws.Cells.Item(2609, 4).Value = 0.00137602447160157
$ws.Cells.Item(2609, 5).Value = 13
$ws.Cells.Item(2609, 6).Value = 0.00264550264550265
$ws.Cells.Item(2609, 7).Value = 0
$ws.Cells.Item(2609, 8).Value = 0

$ws.Cells.Item(2610, 1).Value = 44171
$ws.Cells.Item(2610, 2).Value = '0-10 years'
$ws.Cells.Item(2610, 3).Value = 19909
$ws.Cells.Item(2610, 4).Value = 0.04969869743431
$ws.Cells.Item(2610, 5).Value = 136
$ws.Cells.Item(2610, 6).Value = 0.0442708333333333
$ws.Cells.Item(2610, 7).Value = 4
$ws.Cells.Item(2610, 8).Value = 0

$ws.Cells.Item(2611, 1).Value = 44171
$ws.Cells.Item(2611, 2).Value = '11-20 years'
$ws.Cells.Item(2611, 3).Value = 52714
$ws.Cells.Item(2611, 4).Value = 0.131589589459652
$ws.Cells.Item(2611, 5).Value = 329
$ws.Cells.Item(2611, 6).Value = 0.107096354166667
$ws.Cells.Item(2611, 7).Value = 2
$ws.Cells.Item(2611, 8).Value = 0

$ws.Cells.Item(2612, 1).Value = 44171
$ws.Cells.Item(2612, 2).Value = '21-30 years'
$ws.Cells.Item(2612, 3).Value = 76157
$ws.Cells.Item(2612, 4).Value = 0.190110186373236
$ws.Cells.Item(2612, 5).Value = 536
$ws.Cells.Item(2612, 6).Value = 0.174479166666667
$ws.Cells.Item(2612, 7).Value = 31
$ws.Cells.Item(2612, 8).Value = 0

$ws.Cells.Item(2613, 1).Value = 44171
$ws.Cells.Item(2613, 2).Value = '31-40 years'
$ws.Cells.Item(2613, 3).Value = 62994
$ws.Cells.Item(2613, 4).Value = 0.157251481549898
$ws.Cells.Item(2613, 5).Value = 499
$ws.Cells.Item(2613, 6).Value = 0.162434895833333
$ws.Cells.Item(2613, 7).Value = 62
$ws.Cells.Item(2613, 8).Value = 0

$ws.Cells.Item(2614, 1).Value = 44171
$ws.Cells.Item(2614, 2).Value = '41-50 years'
$ws.Cells.Item(2614, 3).Value = 59785
$ws.Cells.Item(2614, 4).Value = 0.149240877297214
$ws.Cells.Item(2614, 5).Value = 493
$ws.Cells.Item(2614, 6).Value = 0.160481770833333
$ws.Cells.Item(2614, 7).Value = 164
$ws.Cells.Item(2614, 8).Value = 2

$ws.Cells.Item(2615, 1).Value = 44171
$ws.Cells.Item(2615, 2).Value = '51-60 years'
$ws.Cells.Item(2615, 3).Value = 54495
$ws.Cells.Item(2615, 4).Value = 0.136035487301357
$ws.Cells.Item(2615, 5).Value = 446
$ws.Cells.Item(2615, 6).Value = 0.145182291666667
$ws.Cells.Item(2615, 7).Value = 455
$ws.Cells.Item(2615, 8).Value = 4

$ws.Cells.Item(2616, 1).Value = 44171
$ws.Cells.Item(2616, 2).Value = '61-70 years'
$ws.Cells.Item(2616, 3).Value = 38534
$ws.Cells.Item(2616, 4).Value = 0.0961921546503442
$ws.Cells.Item(2616, 5).Value = 359
$ws.Cells.Item(2616, 6).Value = 0.116861979166667
$ws.Cells.Item(2616, 7).Value = 920
$ws.Cells.Item(2616, 8).Value = 4

$ws.Cells.Item(2617, 1).Value = 44171
$ws.Cells.Item(2617, 2).Value = '71-80 years'
$ws.Cells.Item(2617, 3).Value = 23124
$ws.Cells.Item(2617, 4).Value = 0.0577242794450241
$ws.Cells.Item(2617, 5).Value = 175
$ws.Cells.Item(2617, 6).Value = 0.0569661458333333
$ws.Cells.Item(2617, 7).Value = 1489
$ws.Cells.Item(2617, 8).Value = 5

$ws.Cells.Item(2618, 1).Value = 44171
$ws.Cells.Item(2618, 2).Value = '81+ years'
$ws.Cells.Item(2618, 3).Value = 12336
$ws.Cells.Item(2618, 4).Value = 0.0307942705082952
$ws.Cells.Item(2618, 5).Value = 100
$ws.Cells.Item(2618, 6).Value = 0.0325520833333333
$ws.Cells.Item(2618, 7).Value = 1816
$ws.Cells.Item(2618, 8).Value = 23

$ws.Cells.Item(2619, 1).Value = 44171
$ws.Cells.Item(2619, 2).Value = 'Pending'
$ws.Cells.Item(2619, 3).Value = 546
$ws.Cells.Item(2619, 4).Value = 0.00136297598066871
$ws.Cells.Item(2619, 5).Value = -1
$ws.Cells.Item(2619, 6).Value = -0.000325520833333333
$ws.Cells.Item(2619, 7).Value = 0
$ws.Cells.Item(2619, 8).Value = 0

$ws.Cells.Item(2620, 1).Value = 44172
$ws.Cells.Item(2620, 2).Value = '0-10 years'
$ws.Cells.Item(2620, 3).Value = 20367
$ws.Cells.Item(2620, 4).Value = 0.049829961099014
$ws.Cells.Item(2620, 5).Value = 458
$ws.Cells.Item(2620, 6).Value = 0.0562930186823992
$ws.Cells.Item(2620, 7).Value = 4
$ws.Cells.Item(2620, 8).Value = 0

$ws.Cells.Item(2621, 1).Value = 44172
$ws.Cells.Item(2621, 2).Value = '11-20 years'
$ws.Cells.Item(2621, 3).Value = 53725
$ws.Cells.Item(2621, 4).Value = 0.131443740366501
$ws.Cells.Item(2621, 5).Value = 1011
$ws.Cells.Item(2621, 6).Value = 0.124262536873156
$ws.Cells.Item(2621, 7).Value = 2
$ws.Cells.Item(2621, 8).Value = 0

$ws.Cells.Item(2622, 1).Value = 44172
$ws.Cells.Item(2622, 2).Value = '21-30 years'
$ws.Cells.Item(2622, 3).Value = 77496
$ws.Cells.Item(2622, 4).Value = 0.18960193770949
$ws.Cells.Item(2622, 5).Value = 1339
$ws.Cells.Item(2622, 6).Value = 0.164577187807276
$ws.Cells.Item(2622, 7).Value = 31
$ws.Cells.Item(2622, 8).Value = 0

$ws.Cells.Item(2623, 1).Value = 44172
$ws.Cells.Item(2623, 2).Value = '31-40 years'
$ws.Cells.Item(2623, 3).Value = 64188
$ws.Cells.Item(2623, 4).Value = 0.15704254642429
$ws.Cells.Item(2623, 5).Value = 1194
$ws.Cells.Item(2623, 6).Value = 0.146755162241888
$ws.Cells.Item(2623, 7).Value = 63
$ws.Cells.Item(2623, 8).Value = 1

$ws.Cells.Item(2624, 1).Value = 44172
$ws.Cells.Item(2624, 2).Value = '41-50 years'
$ws.Cells.Item(2624, 3).Value = 61083
$ws.Cells.Item(2624, 4).Value = 0.149445844444988
$ws.Cells.Item(2624, 5).Value = 1298
$ws.Cells.Item(2624, 6).Value = 0.159537856440511
$ws.Cells.Item(2624, 7).Value = 166
$ws.Cells.Item(2624, 8).Value = 2

$ws.Cells.Item(2625, 1).Value = 44172
$ws.Cells.Item(2625, 2).Value = '51-60 years'
$ws.Cells.Item(2625, 3).Value = 55686
$ws.Cells.Item(2625, 4).Value = 0.136241528637487
$ws.Cells.Item(2625, 5).Value = 1191
$ws.Cells.Item(2625, 6).Value = 0.146386430678466
$ws.Cells.Item(2625, 7).Value = 457
$ws.Cells.Item(2625, 8).Value = 2

$ws.Cells.Item(2626, 1).Value = 44172
$ws.Cells.Item(2626, 2).Value = '61-70 years'
$ws.Cells.Item(2626, 3).Value = 39448
$ws.Cells.Item(2626, 4).Value = 0.0965135908790644
$ws.Cells.Item(2626, 5).Value = 914
$ws.Cells.Item(2626, 6).Value = 0.112340216322517
$ws.Cells.Item(2626, 7).Value = 932
$ws.Cells.Item(2626, 8).Value = 12

$ws.Cells.Item(2627, 1).Value = 44172
$ws.Cells.Item(2627, 2).Value = '71-80 years'
$ws.Cells.Item(2627, 3).Value = 23604
$ws.Cells.Item(2627, 4).Value = 0.0577496146600445
$ws.Cells.Item(2627, 5).Value = 480
$ws.Cells.Item(2627, 6).Value = 0.0589970501474926
$ws.Cells.Item(2627, 7).Value = 1510
$ws.Cells.Item(2627, 8).Value = 21

$ws.Cells.Item(2628, 1).Value = 44172
$ws.Cells.Item(2628, 2).Value = '81+ years'
$ws.Cells.Item(2628, 3).Value = 12584
$ws.Cells.Item(2628, 4).Value = 0.030788050791476
$ws.Cells.Item(2628, 5).Value = 248
$ws.Cells.Item(2628, 6).Value = 0.0304818092428712
$ws.Cells.Item(2628, 7).Value = 1844
$ws.Cells.Item(2628, 8).Value = 28

$ws.Cells.Item(2629, 1).Value = 44172
$ws.Cells.Item(2629, 2).Value = 'Pending'
$ws.Cells.Item(2629, 3).Value = 549
$ws.Cells.Item(2629, 4).Value = 0.00134318498764466
$ws.Cells.Item(2629, 5).Value = 3
$ws.Cells.Item(2629, 6).Value = 0.000368731563421829
$ws.Cells.Item(2629, 7).Value = 0
$ws.Cells.Item(2629, 8).Value = 0

# Update the workbook-level defined name range and sheet dimension to
# reflect the newly-added rows (A1:H2629).
$wb.Names.Item("ALL_AGE_FINAL").RefersTo = "='ALL_AGE_FINAL'!`$A`$1:`$H`$2629"

Write-Output "done"
